$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: row, B value, C value, D value
$data = @(
    @(3, 68035, 9929, 10908),
    @(4, 32385, 3774, 3830),
    @(5, 113941, 8075, 10466),
    @(6, 2292, 760, 127),
    @(7, 67616, 11539, 10235),
    @(8, 8325, 1732, 1671),
    @(9, 8620, 1550, 1119),
    @(10, 3253, 557, 160),
    @(11, 459, 388, 11),
    @(12, 0, 0, 0),
    @(13, 1485, 453, 299),
    @(14, 4680, 2008, 1475),
    @(15, 9099, 3110, 1344),
    @(16, 5198, 2163, 862),
    @(17, 2959, 1036, 253),
    @(18, 25888, 4009, 4752),
    @(19, 4447, 1155, 900),
    @(20, 32145, 3521, 5554),
    @(21, 667, 497, 30),
    @(22, 26901, 3257, 4288),
    @(23, 1839, 733, 283),
    @(24, 27818, 4785, 4868),
    @(25, 107705, 12021, 12033),
    @(26, 7772, 2830, 1084),
    @(27, 0, 0, 0),
    @(28, 8888, 1737, 1696),
    @(29, 3304, 511, 712),
    @(30, 24294, 4425, 4448),
    @(31, 894, 314, 349),
    @(32, 3716, 2403, 487),
    @(33, 23276, 4708, 4323),
    @(34, 16488, 3984, 4000),
    @(35, 9771, 1175, 2082),
    @(36, 80501, 8744, 8325),
    @(37, 12733, 4061, 2096),
    @(38, 30954, 3014, 4423),
    @(39, 1376, 1202, 260),
    @(40, 2484, 350, 1151),
    @(41, 3857, 495, 232),
    @(42, 15558, 327, 389),
    @(43, 391, 150, 89),
    @(44, 1306, 64, 32),
    @(45, 3155, 244, 73),
    @(46, 4975, 1435, 783),
    @(47, 19011, 4978, 3611),
    @(48, 45488, 4889, 6782),
    @(49, 21987, 5110, 1863),
    @(50, 16589, 1554, 2894),
    @(51, 42947, 4455, 5237),
    @(52, 5994, 701, 1350),
    @(53, 21290, 4408, 3535),
    @(54, 3384, 1131, 1687),
    @(55, 3257, 2098, 308),
    @(56, 4969, 1594, 1574),
    @(57, 20462, 7734, 3980),
    @(58, 24267, 1280, 815),
    @(59, 996366, 150143, 141343)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
}
